$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value = "Integral Adaptive Law와 Persistent Excitation에 관하여"
$ws.Range("E28").Value = "https://ropiens.tistory.com/244"

$ws.Range("D32").Value = "[tensorflow in spark] spark를 이용해 tf model을 분산 처리?!"
$ws.Range("E32").Value = "https://dodonam.tistory.com/484"

$ws.Range("D36").Value = "Diffusion models for Time-series"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/442"
